$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns keep their original text formatting
# so values like "1.000" or "30.869.60" are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '30.869.60'
$ws.Range("E2").Value = '  +2.43%  '

$ws.Range("D3").Value = '2.115.88'
$ws.Range("E3").Value = '  +10.20%  '

$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = '334.62'
$ws.Range("E5").Value = '  +4.74%  '

$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.19%  '

$ws.Range("D7").Value = '0.5351'
$ws.Range("E7").Value = '  +5.52%  '

$ws.Range("D8").Value = '0.4415'
$ws.Range("E8").Value = '  +8.27%  '

$ws.Range("D9").Value = '0.09062'
$ws.Range("E9").Value = '  +8.67%  '

$ws.Range("D10").Value = '46.40'
$ws.Range("E10").Value = '  +10.28%  '

$ws.Range("D11").Value = '1.183'
$ws.Range("E11").Value = '  +5.93%  '

$ws.Range("D12").Value = '25.38'
$ws.Range("E12").Value = '  +4.76%  '

$ws.Range("D13").Value = '2.122.17'
$ws.Range("E13").Value = '  +10.41%  '

$ws.Range("D14").Value = '6.779'
$ws.Range("E14").Value = '  +5.48%  '

$ws.Range("D15").Value = '7.832'
$ws.Range("E15").Value = '  +7.93%  '

$ws.Range("D16").Value = '98.11'
$ws.Range("E16").Value = '  +5.84%  '

$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  -0.30%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = '0.00001139'
$ws.Range("E18").Value = '  +3.93%  '

$ws.Range("D19").Value = '0.06656'
$ws.Range("E19").Value = '  +2.11%  '

$ws.Range("D20").Value = '19.23'
$ws.Range("E20").Value = '  +3.92%  '

$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '0.9996'
$ws.Range("E21").Value = '  -0.22%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '6.394'
$ws.Range("E22").Value = '  +7.29%  '

$ws.Range("D23").Value = '30.975.20'
$ws.Range("E23").Value = '  +2.76%  '

$ws.Range("D24").Value = '12.20'
$ws.Range("E24").Value = '  +7.36%  '

$ws.Range("D25").Value = '2.367.59'
$ws.Range("E25").Value = '  +10.55%  '

$ws.Range("D26").Value = '2.266'
$ws.Range("E26").Value = '  +3.17%  '

$ws.Range("D27").Value = '22.90'
$ws.Range("E27").Value = '  +4.56%  '

$ws.Range("D28").Value = '2.589'
$ws.Range("E28").Value = '  +14.43%  '

$ws.Range("D29").Value = '163.84'
$ws.Range("E29").Value = '  +0.67%  '

$ws.Range("D30").Value = '134.07'
$ws.Range("E30").Value = '  +4.07%  '

$ws.Range("D31").Value = '1.175'
$ws.Range("E31").Value = '  +3.05%  '

$ws.Range("D32").Value = '0.1082'
$ws.Range("E32").Value = '  +3.25%  '

$ws.Range("D33").Value = '6.272'
$ws.Range("E33").Value = '  +5.31%  '

$ws.Range("D34").Value = '4.004'
$ws.Range("E34").Value = '  +5.68%  '

$ws.Range("D35").Value = '1.535'
$ws.Range("E35").Value = '  +26.38%  '

$ws.Range("D36").Value = '0.02622'
$ws.Range("E36").Value = '  +7.04%  '

$ws.Range("D37").Value = '13.33'
$ws.Range("E37").Value = '  +16.18%  '

$ws.Range("D38").Value = '5.589'
$ws.Range("E38").Value = '  +5.26%  '

$ws.Range("D39").Value = '9.615'
$ws.Range("E39").Value = '  +11.95%  '

$ws.Range("D40").Value = '0.06755'
$ws.Range("E40").Value = '  +4.82%  '

$ws.Range("D41").Value = '0.2281'
$ws.Range("E41").Value = '  +6.27%  '

$ws.Range("D42").Value = '0.6880'

$ws.Range("D43").Value = '1.258'
$ws.Range("E43").Value = '  +3.75%  '

$ws.Range("D44").Value = '14.22'
$ws.Range("E44").Value = '  +6.00%  '

$ws.Range("D45").Value = '0.6459'
$ws.Range("E45").Value = '  +6.75%  '

$ws.Range("D46").Value = '0.9997'
$ws.Range("E46").Value = '  -0.08%  '

$ws.Range("E47").Value = '  +3.47%  '

$ws.Range("D48").Value = '3.684'
$ws.Range("E48").Value = '  +1.66%  '

$ws.Range("D49").Value = '1.288'
$ws.Range("E49").Value = '  +6.37%  '

$ws.Range("D50").Value = '83.35'
$ws.Range("E50").Value = '  +6.26%  '

$ws.Range("D51").Value = '1.176'
$ws.Range("E51").Value = '  +3.59%  '
